# 自动更新Excel文件
# 每日刷新"剩余"天数(列E)：剩余 = 总天(列D) - (当前日期 - 开始时间(列F))
# 当剩余归零(即前一日剩余=1)时，说明周期结束，重新开始新周期：
#   开始时间(列F) 重置为当前日期，剩余(列E) 重置为总天(列D)
# 本次更新将"当前日期"从 2025-10-26 推进到 2025-10-27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; E=7; F=20251020},
    @{Row=3; E=7; F=20251020},
    @{Row=4; E=7; F=20251020},
    @{Row=5; E=9; F=20251026},
    @{Row=6; E=7; F=20251020},
    @{Row=7; E=9; F=20251026},
    @{Row=8; E=7; F=20251020},
    @{Row=9; E=9; F=20251026},
    @{Row=10; E=7; F=20251027},
    @{Row=11; E=7; F=20251020},
    @{Row=12; E=9; F=20251026},
    @{Row=13; E=7; F=20251020},
    @{Row=14; E=7; F=20251020},
    @{Row=15; E=7; F=20251020},
    @{Row=16; E=3; F=20251020},
    @{Row=17; E=9; F=20251026},
    @{Row=18; E=2; F=20251019},
    @{Row=19; E=2; F=20251019},
    @{Row=20; E=2; F=20251019},
    @{Row=21; E=2; F=20251019},
    @{Row=22; E=9; F=20251026},
    @{Row=23; E=9; F=20251026},
    @{Row=24; E=9; F=20251026},
    @{Row=25; E=9; F=20251026},
    @{Row=26; E=9; F=20251026},
    @{Row=27; E=1; F=20251021},
    @{Row=28; E=2; F=20251019},
    @{Row=29; E=2; F=20251019},
    @{Row=30; E=2; F=20251019},
    @{Row=31; E=2; F=20251019},
    @{Row=32; E=2; F=20251019},
    @{Row=33; E=2; F=20251019},
    @{Row=34; E=2; F=20251019},
    @{Row=35; E=2; F=20251019},
    @{Row=37; E=2; F=20251019},
    @{Row=38; E=2; F=20251019},
    @{Row=39; E=2; F=20251019},
    @{Row=40; E=7; F=20251027},
    @{Row=41; E=7; F=20251027},
    @{Row=42; E=2; F=20251019},
    @{Row=43; E=9; F=20251026},
    @{Row=44; E=7; F=20251027},
    @{Row=45; E=9; F=20251026},
    @{Row=46; E=7; F=20251027},
    @{Row=47; E=2; F=20251019},
    @{Row=48; E=7; F=20251027},
    @{Row=49; E=1; F=20251021},
    @{Row=50; E=7; F=20251024},
    @{Row=51; E=7; F=20251024},
    @{Row=52; E=7; F=20251024},
    @{Row=53; E=7; F=20251024},
    @{Row=54; E=7; F=20251024},
    @{Row=55; E=7; F=20251024},
    @{Row=56; E=7; F=20251024},
    @{Row=57; E=7; F=20251024},
    @{Row=58; E=1; F=20251018},
    @{Row=59; E=1; F=20251018},
    @{Row=60; E=1; F=20251018},
    @{Row=61; E=1; F=20251021},
    @{Row=62; E=1; F=20251018},
    @{Row=63; E=1; F=20251018},
    @{Row=64; E=1; F=20251018},
    @{Row=65; E=2; F=20251019},
    @{Row=66; E=2; F=20251019},
    @{Row=67; E=2; F=20251019},
    @{Row=68; E=2; F=20251019},
    @{Row=69; E=2; F=20251019},
    @{Row=70; E=3; F=20251020},
    @{Row=71; E=3; F=20251020},
    @{Row=72; E=3; F=20251020},
    @{Row=73; E=3; F=20251020},
    @{Row=74; E=3; F=20251020},
    @{Row=75; E=3; F=20251020},
    @{Row=76; E=3; F=20251020},
    @{Row=77; E=6; F=20251023},
    @{Row=78; E=6; F=20251023},
    @{Row=79; E=6; F=20251023},
    @{Row=80; E=6; F=20251023},
    @{Row=81; E=6; F=20251023},
    @{Row=82; E=6; F=20251023},
    @{Row=83; E=6; F=20251023},
    @{Row=84; E=6; F=20251023},
    @{Row=85; E=6; F=20251023},
    @{Row=86; E=6; F=20251023},
    @{Row=87; E=7; F=20251027},
    @{Row=88; E=7; F=20251027},
    @{Row=89; E=7; F=20251027},
    @{Row=90; E=7; F=20251027},
    @{Row=91; E=9; F=20251026},
    @{Row=92; E=7; F=20251027},
    @{Row=93; E=6; F=20251023},
    @{Row=94; E=3; F=20251023},
    @{Row=95; E=5; F=20251022},
    @{Row=96; E=3; F=20251020},
    @{Row=97; E=3; F=20251020},
    @{Row=98; E=3; F=20251020},
    @{Row=99; E=3; F=20251020}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = $u.F
}
